$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Split the old "Terms Typically Offered" column (D) requirement info into
# three new columns: Corequisites, Concurrent, Recommended. Insert three
# fresh columns before D so the old D column (and its data) shifts right
# to G, then populate the new D:F columns.
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("D1").EntireColumn.Insert()
$ws.Range("D1").EntireColumn.Insert()

# New header row labels
$ws.Range("D1").Value = "Corequisites"
$ws.Range("E1").Value = "Concurrent"
$ws.Range("F1").Value = "Recommended"

# New data cells for each course row (2-5) default to "NA"
$ws.Range("D2:F5").Value = "NA"
